$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate previously-empty rows 51-53 with new data
$ws.Range("A51").Value = 46013
$ws.Range("B51").Value = 573
$ws.Range("C51").Value = 10
$ws.Range("D51").Value = 563

$ws.Range("A52").Value = 46010
$ws.Range("B52").Value = 746
$ws.Range("C52").Value = 17
$ws.Range("D52").Value = 729

$ws.Range("A53").Value = 46014
$ws.Range("B53").Value = 656
$ws.Range("C53").Value = 20
$ws.Range("D53").Value = 636

# Update the view state to match the scrolled/selected position after edit
$ws.Application.ActiveWindow.ScrollRow = 44
$ws.Range("A53:D53").Select()
